$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: write the new values as text (apostrophe-prefixed so Excel
# keeps the numeric-looking strings as literal text, matching the
# original inlineStr cell type instead of silently parsing them as numbers).
$ws.Range("D2").Value = "'306.33"
$ws.Range("E2").Value = "'0.37%"
$ws.Range("D3").Value = "'36.15"
$ws.Range("E3").Value = "'0.91%"
$ws.Range("D4").Value = "'5.059"
$ws.Range("E4").Value = "'-0.13%"
$ws.Range("D5").Value = "'0.08049"
$ws.Range("E5").Value = "'-0.07%"
$ws.Range("D6").Value = "'1.876"
$ws.Range("E6").Value = "'-2.13%"
$ws.Range("D7").Value = "'7.792"
$ws.Range("E7").Value = "'-0.58%"
$ws.Range("D8").Value = "'0.9264"
$ws.Range("E8").Value = "'-0.38%"
$ws.Range("D9").Value = "'0.1440"
$ws.Range("E9").Value = "'8.52%"
$ws.Range("D10").Value = "'0.1911"
$ws.Range("E10").Value = "'0.43%"
$ws.Range("D11").Value = "'0.09023"
$ws.Range("E11").Value = "'-1.78%"
$ws.Range("D12").Value = "'0.03441"
$ws.Range("E12").Value = "'-0.86%"
$ws.Range("D13").Value = "'0.09909"
$ws.Range("E13").Value = "'-0.09%"
$ws.Range("D14").Value = "'0.001404"
$ws.Range("E14").Value = "'-0.99%"
$ws.Range("D15").Value = "'0.006035"
$ws.Range("E15").Value = "'-9.00%"
$ws.Range("D16").Value = "'3.841"
$ws.Range("E16").Value = "'6.35%"
$ws.Range("D17").Value = "'4.126"
$ws.Range("E17").Value = "'-0.85%"
$ws.Range("D18").Value = "'3.392"
$ws.Range("E18").Value = "'13.45%"
$ws.Range("D19").Value = "'0.3448"
$ws.Range("E19").Value = "'0.73%"
$ws.Range("D20").Value = "'0.1335"
$ws.Range("E20").Value = "'-0.13%"
$ws.Range("D21").Value = "'4.822"
$ws.Range("E21").Value = "'-7.08%"
$ws.Range("D23").Value = "'0.04374"
$ws.Range("E23").Value = "'-1.06%"
$ws.Range("E24").Value = "'-0.70%"
$ws.Range("D25").Value = "'0.004295"
$ws.Range("E25").Value = "'-8.73%"
$ws.Range("D27").Value = "'0.0001299"
$ws.Range("E27").Value = "'-0.29%"
$ws.Range("E28").Value = "'41.90%"
$ws.Range("D39").Value = "'0.02004"
$ws.Range("E39").Value = "'0.05%"
$ws.Range("D40").Value = "'0.05113"
$ws.Range("E40").Value = "'-1.62%"
$ws.Range("D41").Value = "'0.007541"
$ws.Range("E41").Value = "'-1.32%"
$ws.Range("D42").Value = "'0.01011"
$ws.Range("E42").Value = "'-0.14%"
$ws.Range("D43").Value = "'0.1359"
$ws.Range("E43").Value = "'-0.38%"
$ws.Range("D44").Value = "'0.002168"
$ws.Range("E44").Value = "'3.04%"
$ws.Range("D45").Value = "'0.009632"
$ws.Range("E45").Value = "'-10.10%"
$ws.Range("D46").Value = "'0.00006243"
$ws.Range("E46").Value = "'-1.00%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.19%"
$ws.Range("D48").Value = "'64.92"
$ws.Range("E48").Value = "'-0.46%"
$ws.Range("D49").Value = "'0.001250"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.19%"

# Step 2: the apostrophe prefix above makes Excel mark the cells with a
# "quote prefix" style. Restore each touched cell to the sheet's plain
# (unstyled) look by pasting the formatting from the neighboring, always-
# plain F:G columns of the same rows - this mirrors the original file,
# where these data cells carry no cell style at all.
$ws.Range("F2:G21").Copy() | Out-Null
$ws.Range("D2:E21").PasteSpecial(-4122) | Out-Null
$ws.Range("F23:G25").Copy() | Out-Null
$ws.Range("D23:E25").PasteSpecial(-4122) | Out-Null
$ws.Range("F27:G28").Copy() | Out-Null
$ws.Range("D27:E28").PasteSpecial(-4122) | Out-Null
$ws.Range("F39:G51").Copy() | Out-Null
$ws.Range("D39:E51").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
